$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.98"
$ws.Range("D3").Value = "'25.01"
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D4").Value = "'5.104"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D5").Value = "'0.05635"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'6.532"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'2.982"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8122"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.8393"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1337"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.06951"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03231"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02842"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09401"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001516"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0005952"
$ws.Range("E16").Value = "15OneONEWorstin24h"
$ws.Range("D17").Value = "'0.006139"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.499"
$ws.Range("E18").Value = "17LEOLEO"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.107"
$ws.Range("E19").Value = "18BTSETokenBTSE"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3171"
$ws.Range("E20").Value = "19BitpandaEcosystemTokenBEST"
$ws.Range("D21").Value = "'0.1292"
$ws.Range("D22").Value = "'3.746"
$ws.Range("D23").Value = "'0.04707"
$ws.Range("D24").Value = "'0.1370"
$ws.Range("D25").Value = "'0.001243"
$ws.Range("D26").Value = "'0.004524"
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("D28").Value = "'0.0001546"
$ws.Range("D40").Value = "'0.03625"
$ws.Range("D41").Value = "'0.006239"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1051"
$ws.Range("D43").Value = "'0.002729"
$ws.Range("D44").Value = "'0.008304"
$ws.Range("D45").Value = "'0.00005302"
$ws.Range("D47").Value = "'0.1801"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Write-Host "Applied 84 cell updates"